$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2005730659025788
$ws.Range("C2").Value = 0.5329512893982808
$ws.Range("J2").Value = 0.01432664756446991
$ws.Range("P2").Value = 0.169054441260745
$ws.Range("S2").Value = 0.0830945558739255
$ws.Range("B3").Value = 0.005208333333333333
$ws.Range("C3").Value = 0.03125
$ws.Range("J3").Value = 0.06770833333333333
$ws.Range("P3").Value = 0.7135416666666666
$ws.Range("S3").Value = 0.1822916666666667
$ws.Range("J4").Value = 0.07547169811320754
$ws.Range("P4").Value = 0.6792452830188679
$ws.Range("S4").Value = 0.2452830188679245
$ws.Range("B6").Value = 0.09523809523809523
$ws.Range("D6").Value = 0.02380952380952381
$ws.Range("F6").Value = 0.04761904761904762
$ws.Range("J6").Value = 0.2857142857142857
$ws.Range("O6").Value = 0.01904761904761905
$ws.Range("Q6").Value = 0.119047619047619
$ws.Range("R6").Value = 0.06666666666666667
$ws.Range("S6").Value = 0.3428571428571429
$ws.Range("B7").Value = 0.1100917431192661
$ws.Range("D7").Value = 0.02293577981651376
$ws.Range("F7").Value = 0.07339449541284404
$ws.Range("J7").Value = 0.1192660550458716
$ws.Range("O7").Value = 0.01376146788990826
$ws.Range("Q7").Value = 0.1834862385321101
$ws.Range("R7").Value = 0.1009174311926606
$ws.Range("S7").Value = 0.3761467889908257
$ws.Range("B8").Value = 0.1457286432160804
$ws.Range("D8").Value = 0.02512562814070352
$ws.Range("F8").Value = 0.03266331658291458
$ws.Range("J8").Value = 0.1080402010050251
$ws.Range("O8").Value = 0.02010050251256281
$ws.Range("Q8").Value = 0.1909547738693467
$ws.Range("R8").Value = 0.1155778894472362
$ws.Range("S8").Value = 0.3618090452261307
$ws.Range("B9").Value = 0.1
$ws.Range("D9").Value = 0.02
$ws.Range("F9").Value = 0.056
$ws.Range("J9").Value = 0.12
$ws.Range("O9").Value = 0.02
$ws.Range("Q9").Value = 0.172
$ws.Range("R9").Value = 0.108
$ws.Range("S9").Value = 0.404
$ws.Range("B10").Value = 0.1113636363636364
$ws.Range("D10").Value = 0.02272727272727273
$ws.Range("E10").Value = 0.002272727272727273
$ws.Range("F10").Value = 0.06515151515151515
$ws.Range("J10").Value = 0.121969696969697
$ws.Range("O10").Value = 0.007575757575757576
$ws.Range("Q10").Value = 0.221969696969697
$ws.Range("R10").Value = 0.0946969696969697
$ws.Range("S10").Value = 0.3522727272727273
$ws.Range("G11").Value = 0.1578947368421053
$ws.Range("J11").Value = 0.07309941520467836
$ws.Range("K11").Value = 0.1929824561403509
$ws.Range("L11").Value = 0.5526315789473685
$ws.Range("S11").Value = 0.02339181286549707
$ws.Range("F12").Value = 0.005181347150259068
$ws.Range("G12").Value = 0.7305699481865285
$ws.Range("J12").Value = 0.2020725388601036
$ws.Range("K12").Value = 0.0155440414507772
$ws.Range("L12").Value = 0.02072538860103627
$ws.Range("S12").Value = 0.02590673575129534
$ws.Range("G13").Value = 0.5957446808510638
$ws.Range("J13").Value = 0.2978723404255319
$ws.Range("S13").Value = 0.1063829787234043
$ws.Range("F15").Value = 0.02926829268292683
$ws.Range("H15").Value = 0.07804878048780488
$ws.Range("I15").Value = 0.1073170731707317
$ws.Range("J15").Value = 0.3560975609756097
$ws.Range("K15").Value = 0.07804878048780488
$ws.Range("M15").Value = 0.02439024390243903
$ws.Range("O15").Value = 0.07804878048780488
$ws.Range("S15").Value = 0.248780487804878
$ws.Range("F16").Value = 0.03153153153153153
$ws.Range("H16").Value = 0.1306306306306306
$ws.Range("I16").Value = 0.1126126126126126
$ws.Range("J16").Value = 0.3288288288288289
$ws.Range("K16").Value = 0.1396396396396396
$ws.Range("M16").Value = 0.02252252252252252
$ws.Range("O16").Value = 0.05855855855855856
$ws.Range("S16").Value = 0.1756756756756757
$ws.Range("F17").Value = 0.01483050847457627
$ws.Range("H17").Value = 0.1694915254237288
$ws.Range("I17").Value = 0.09957627118644068
$ws.Range("J17").Value = 0.4279661016949153
$ws.Range("K17").Value = 0.08050847457627118
$ws.Range("M17").Value = 0.01483050847457627
$ws.Range("N17").Value = 0.00211864406779661
$ws.Range("O17").Value = 0.05084745762711865
$ws.Range("S17").Value = 0.1398305084745763
$ws.Range("F18").Value = 0.00425531914893617
$ws.Range("H18").Value = 0.174468085106383
$ws.Range("I18").Value = 0.1063829787234043
$ws.Range("J18").Value = 0.4170212765957447
$ws.Range("K18").Value = 0.1191489361702128
$ws.Range("M18").Value = 0.02127659574468085
$ws.Range("O18").Value = 0.05531914893617021
$ws.Range("S18").Value = 0.1021276595744681
$ws.Range("F19").Value = 0.01822323462414579
$ws.Range("H19").Value = 0.1738800303720577
$ws.Range("I19").Value = 0.1009870918754746
$ws.Range("J19").Value = 0.3568716780561883
$ws.Range("K19").Value = 0.1207289293849658
$ws.Range("M19").Value = 0.02126044039483675
$ws.Range("N19").Value = 0.0007593014426727411
$ws.Range("O19").Value = 0.06150341685649203
$ws.Range("S19").Value = 0.1457858769931663
